$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted as row 197 (pushing the existing
# rows 197-263 down to 198-264). Insert a blank row at position 197 first,
# preserving/shifting all data below it.
$ws.Rows.Item(197).Insert()

# Populate the newly inserted row 197 with the new record's data.
$ws.Cells.Item(197, 1).Value = 3
$ws.Cells.Item(197, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(197, 3).Value = "Coquimbo"
$ws.Cells.Item(197, 4).Value = 44559
$ws.Cells.Item(197, 5).Value = 5
$ws.Cells.Item(197, 6).Value = 100112043
$ws.Cells.Item(197, 7).Value = "Pepino ensalada"
$ws.Cells.Item(197, 8).Value = "Sin especificar"
$ws.Cells.Item(197, 9).Value = "Primera"
$ws.Cells.Item(197, 10).Value = 114
$ws.Cells.Item(197, 11).Value = 6500
$ws.Cells.Item(197, 12).Value = 7000
$ws.Cells.Item(197, 13).Value = 6746
$ws.Cells.Item(197, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(197, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(197, 16).Value = 96
$ws.Cells.Item(197, 17).Value = 70
$ws.Cells.Item(197, 18).Value = "Hortaliza"
